$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 is newly added - fill in the constant columns (A,B,C,E,F,G,H,I,J,K) that match all other data rows
$ws.Range("A14").Value = 9
$ws.Range("B14").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C14").Value = "Metropolitana"
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100101
$ws.Range("H14").Value = "Berries"
$ws.Range("I14").Value = 100101006
$ws.Range("J14").Value = "Breva"
$ws.Range("K14").Value = "Sin especificar"

# Update D,L,M,N,O,P,Q,R,S,T for rows 3 through 14 per the diff (row 3 is brand-new data; rows 4-14 shift down from the previous row's old values)

# Row 3
$ws.Range("D3").Value = 45274
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 15
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 24000
$ws.Range("Q3").Value = "$/bandeja 5 kilos"
$ws.Range("R3").Value = "Provincia de Huasco"
$ws.Range("S3").Value = 4800
$ws.Range("T3").Value = 5

# Row 4
$ws.Range("D4").Value = 44204
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 110
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7500
$ws.Range("P4").Value = 7318
$ws.Range("Q4").Value = "$/bandeja 7 kilos"
$ws.Range("R4").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S4").Value = 1045
$ws.Range("T4").Value = 7

# Row 5
$ws.Range("D5").Value = 44187
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 45
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 14000
$ws.Range("Q5").Value = "$/bandeja 7 kilos"
$ws.Range("R5").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S5").Value = 2000
$ws.Range("T5").Value = 7

# Row 6
$ws.Range("D6").Value = 44187
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = "$/bandeja 7 kilos"
$ws.Range("R6").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S6").Value = 1714
$ws.Range("T6").Value = 7

# Row 7
$ws.Range("D7").Value = 44553
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 22000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 22000
$ws.Range("Q7").Value = "$/bandeja 6 kilos"
$ws.Range("R7").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S7").Value = 3667
$ws.Range("T7").Value = 6

# Row 8
$ws.Range("D8").Value = 44553
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 150
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("Q8").Value = "$/bandeja 6 kilos"
$ws.Range("R8").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S8").Value = 3000
$ws.Range("T8").Value = 6

# Row 9
$ws.Range("D9").Value = 44572
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 65
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("Q9").Value = "$/bandeja 6 kilos"
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 3333
$ws.Range("T9").Value = 6

# Row 10
$ws.Range("D10").Value = 44558
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 22000
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 22000
$ws.Range("Q10").Value = "$/bandeja 6 kilos"
$ws.Range("R10").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S10").Value = 3667
$ws.Range("T10").Value = 6

# Row 11
$ws.Range("D11").Value = 44558
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 25
$ws.Range("N11").Value = 18000
$ws.Range("O11").Value = 18000
$ws.Range("P11").Value = 18000
$ws.Range("Q11").Value = "$/bandeja 6 kilos"
$ws.Range("R11").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S11").Value = 3000
$ws.Range("T11").Value = 6

# Row 12
$ws.Range("D12").Value = 44550
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 24000
$ws.Range("O12").Value = 24000
$ws.Range("P12").Value = 24000
$ws.Range("Q12").Value = "$/bandeja 7 kilos"
$ws.Range("R12").Value = "Región Metropolitana"
$ws.Range("S12").Value = 3429
$ws.Range("T12").Value = 7

# Row 13
$ws.Range("D13").Value = 44189
$ws.Range("L13").Value = "Especial"
$ws.Range("M13").Value = 20
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range("Q13").Value = "$/bandeja 7 kilos"
$ws.Range("R13").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S13").Value = 2143
$ws.Range("T13").Value = 7

# Row 14
$ws.Range("D14").Value = 44189
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 30
$ws.Range("N14").Value = 13000
$ws.Range("O14").Value = 13000
$ws.Range("P14").Value = 13000
$ws.Range("Q14").Value = "$/bandeja 7 kilos"
$ws.Range("R14").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S14").Value = 1857
$ws.Range("T14").Value = 7
